# DPLKKPS002-009 .. DPLKKPS002-011 and DPLKKPS002-013
# Point the environment's URL text at the new host (192.168.168.107) and
# tidy up the now-unused border/fill formatting on the browser-icon cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the displayed URL text in cell D2 (the hyperlink itself is left
# pointing at its original target; only the visible text changes).
$ws.Range("D2").Value = "http://192.168.168.107/"

# Clear the leftover border/fill formatting on E2 (BROWSER_ICONS value),
# leaving it unlocked (Protection: Locked = False).
$ws.Range("E2").Borders.LineStyle = -4142
$ws.Range("E2").Interior.Pattern = -4142

# Move the active selection to F2, matching the saved view state.
$ws.Range("F2").Select()
